$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws 'D2' '28.893.43'
$ws.Range('E2').Value = '  -2.78%  '
Set-TextValue $ws 'D3' '1.889.61'
$ws.Range('E3').Value = '  -6.00%  '
$ws.Range('E4').Value = '  -0.89%  '
Set-TextValue $ws 'D5' '323.57'
$ws.Range('E5').Value = '  -2.08%  '
Set-TextValue $ws 'D6' '1.001'
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  -2.80%  '
Set-TextValue $ws 'D8' '0.3812'
$ws.Range('E8').Value = '  -4.22%  '
$ws.Range('E9').Value = '  -3.34%  '
Set-TextValue $ws 'D10' '0.07721'
$ws.Range('E10').Value = '  -3.54%  '
Set-TextValue $ws 'D11' '0.9658'
$ws.Range('E11').Value = '  -4.57%  '
Set-TextValue $ws 'D12' '22.04'
$ws.Range('E12').Value = '  -3.50%  '
Set-TextValue $ws 'D13' '1.883.49'
$ws.Range('E13').Value = '  -8.91%  '
Set-TextValue $ws 'D14' '6.942'
$ws.Range('E14').Value = '  -4.75%  '
Set-TextValue $ws 'D15' '5.668'
$ws.Range('E15').Value = '  -4.26%  '
Set-TextValue $ws 'D16' '0.07054'
$ws.Range('E16').Value = '  -1.84%  '
Set-TextValue $ws 'D17' '1.003'
$ws.Range('E17').Value = '  -0.73%  '
Set-TextValue $ws 'D18' '83.34'
$ws.Range('E18').Value = '  -6.88%  '
Set-TextValue $ws 'D19' '0.000009510'
$ws.Range('E19').Value = '  -5.28%  '
Set-TextValue $ws 'D20' '16.67'
$ws.Range('E20').Value = '  -4.90%  '
Set-TextValue $ws 'D21' '1.001'
$ws.Range('E21').Value = '  -0.79%  '
Set-TextValue $ws 'D22' '28.840.34'
$ws.Range('E22').Value = '  -3.24%  '
Set-TextValue $ws 'D23' '5.291'
$ws.Range('E23').Value = '  -5.15%  '
Set-TextValue $ws 'D24' '10.89'
$ws.Range('E24').Value = '  -3.91%  '
Set-TextValue $ws 'D25' '2.123.33'
$ws.Range('E25').Value = '  -7.38%  '
Set-TextValue $ws 'D26' '2.079'
$ws.Range('E26').Value = '  -3.11%  '
Set-TextValue $ws 'D27' '156.35'
$ws.Range('E27').Value = '  -2.06%  '
Set-TextValue $ws 'D28' '19.04'
$ws.Range('E28').Value = '  -3.94%  '
Set-TextValue $ws 'D29' '5.580'
$ws.Range('E29').Value = '  -7.33%  '
Set-TextValue $ws 'D30' '117.08'
$ws.Range('E30').Value = '  -3.50%  '
Set-TextValue $ws 'D31' '1.811'
$ws.Range('E31').Value = '  -8.38%  '
Set-TextValue $ws 'D32' '0.09269'
$ws.Range('E32').Value = '  -2.48%  '
Set-TextValue $ws 'D33' '0.8500'
$ws.Range('E33').Value = '  -5.92%  '
Set-TextValue $ws 'D34' '5.069'
$ws.Range('E34').Value = '  -4.90%  '
Set-TextValue $ws 'D35' '1.231'
$ws.Range('E35').Value = '  -8.95%  '
$ws.Range('E36').Value = '  -6.02%  '
Set-TextValue $ws 'D37' '0.05679'
$ws.Range('E37').Value = '  -3.31%  '
Set-TextValue $ws 'D38' '1.143'
$ws.Range('E38').Value = '  -3.64%  '
$ws.Range('E39').Value = '  -0.86%  '
Set-TextValue $ws 'D40' '0.02034'
$ws.Range('E40').Value = '  -5.16%  '
Set-TextValue $ws 'D41' '7.401'
$ws.Range('E41').Value = '  -7.22%  '
$ws.Range('E42').Value = '  -5.70%  '
$ws.Range('E43').Value = '  -4.71%  '
Set-TextValue $ws 'D44' '0.000002886'
$ws.Range('E44').Value = '  -11.26%  '
Set-TextValue $ws 'D45' '9.187'
$ws.Range('E45').Value = '  -7.61%  '
Set-TextValue $ws 'D46' '2.696'
$ws.Range('E46').Value = '  +0.80%  '
Set-TextValue $ws 'D47' '0.5166'
$ws.Range('E47').Value = '  -4.77%  '
Set-TextValue $ws 'D48' '11.23'
$ws.Range('E48').Value = '  -8.31%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D49' '2.077'
$ws.Range('E49').Value = '  -4.52%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D50' '0.06794'
$ws.Range('E50').Value = '  -3.28%  '
Set-TextValue $ws 'D51' '111.19'
$ws.Range('E51').Value = '  -3.41%  '
